$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark (it gets re-inserted
#     at its new location further down in this same edit) ---
$d.Bookmarks.Item("_GoBack").Delete()

# --- Step 2: locate "emulator." at the end of the underlined heading
#     run, so we don't depend on any hard-coded character offsets ---
$found = $d.Content
$found.Find.ClearFormatting()
$found.Find.Execute("emulator.", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
$periodStart = $found.End - 1   # position of the trailing "."
$periodEnd = $found.End

# split off the trailing period from the run
$periodRange = $d.Range($periodStart, $periodEnd)
$periodRange.Text = ""

# --- Step 3: insert the new text " (for demo purposes only)." right
#     after "...emulator" (replacing the period we just removed) ---
$insPoint = $d.Range($periodStart, $periodStart)
$insPoint.InsertAfter(" (for demo purposes only).")

# --- Step 4: re-insert the "_GoBack" bookmark between the new
#     "(for demo purposes only)" text and the final "." ---
$bmPos = $periodStart + " (for demo purposes only)".Length
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

# --- Step 5: match formatting (14pt, single underline) of the
#     surrounding heading text for all of the newly inserted text ---
$fullLen = " (for demo purposes only).".Length
$formatRange = $d.Range($periodStart, $periodStart + $fullLen)
$formatRange.Font.Size = 14
$formatRange.Font.Underline = 1

# --- Step 6: re-set the final "." run's text in place so the writer
#     drops the (unneeded) xml:space="preserve" it inherited from the
#     " (for demo purposes only)." insertion ---
$finalPeriod = $d.Range($bmPos, $bmPos + 1)
$finalPeriod.Text = "."
